$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for c8dcfb19...md on Overview,
# shared with "Correspond Handoff Datetime" for c8dcfb19...md on de-de.
$wsOverview.Range("G2").Value = "2016-11-08 23:22:32"
$wsDeDe.Range("H2").Value = "2016-11-08 23:22:32"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for c8dcfb19...md
$wsZhCn.Range("H2").Value = "2016-11-08 23:22:18"
$wsZhCn.Range("K2").Value = "2016-11-08 23:23:15"

# de-de: Correspond Handback DateTime for c8dcfb19...md
$wsDeDe.Range("K2").Value = "2016-11-08 23:23:34"
